$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Panel count" (B) to match Max (D) and clear the now-resolved
# "Missing numbers" (E) text for each chapter row whose missing panels
# were found (per commit: Chapter 3 panels 16, 25, 31 added, etc.).
$updates = @{
    54 = 24
    61 = 20
    62 = 17
    63 = 21
    64 = 22
    65 = 20
    66 = 20
    67 = 20
    68 = 20
    69 = 20
    70 = 18
    71 = 22
    72 = 21
    73 = 24
    74 = 23
    75 = 20
    76 = 20
    77 = 22
    78 = 25
    79 = 21
    80 = 21
    81 = 22
    82 = 20
    83 = 22
    84 = 20
    85 = 20
    86 = 20
    87 = 19
    88 = 19
    89 = 20
    90 = 22
    91 = 25
    92 = 20
    93 = 21
    95 = 21
    97 = 21
    99 = 23
    100 = 23
    102 = 24
    103 = 24
    104 = 20
    105 = 22
    107 = 22
    108 = 19
    109 = 20
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
    $ws.Range("E$row").Value = ""
}
